$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "Davangere"
$ws.Range("G15").Value = "Davangere"
$ws.Range("G17").Value = "Davangere"
$ws.Range("G18").Value = "Vijayapura (Bijapur)"
$ws.Range("G23").Value = "Vijayapura (Bijapur)"
$ws.Range("G30").Value = "Davangere"
$ws.Range("G39").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G42").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G43").Value = "Davangere"
$ws.Range("G45").Value = "Kalaburagi (Gulbarga)"

$ws.Range("F21").ClearContents()
